$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 0.029424
$ws.Range("H2").Value = 0.08827199999999999
$ws.Range("I2").Value = 0.1473063425232919
$ws.Range("J2").Value = 0.1473063425232919
$ws.Range("M2").Value = 1.363346333333333
$ws.Range("N2").Value = 4.090039
$ws.Range("O2").Value = 0.02430403345239443
$ws.Range("P2").Value = 0.02430403345239443
$ws.Range("Q2").Value = 0.04011510251199999
$ws.Range("R2").Value = 0.3610359226079999
$ws.Range("S2").Value = 0.003580138276435959
$ws.Range("T2").Value = 0.00358013827643596
$ws.Range("G3").Value = 0.029424
$ws.Range("H3").Value = 0.08827199999999999
$ws.Range("I3").Value = 0.1473063425232919
$ws.Range("J3").Value = 0.1473063425232919
$ws.Range("O3").Value = 0.679596855668023
$ws.Range("P3").Value = 0.679596855668023
$ws.Range("Q3").Value = 1.121710829824
$ws.Range("R3").Value = 10.095397468416
$ws.Range("S3").Value = 0.100108927198786
$ws.Range("T3").Value = 0.100108927198786
$ws.Range("G4").Value = 0.029424
$ws.Range("H4").Value = 0.08827199999999999
$ws.Range("I4").Value = 0.1473063425232919
$ws.Range("J4").Value = 0.1473063425232919
$ws.Range("N4").Value = 49.82946200000001
$ws.Range("O4").Value = 0.2960991108795826
$ws.Range("P4").Value = 0.2960991108795826
$ws.Range("Q4").Value = 0.488727363296
$ws.Range("R4").Value = 4.398546269664
$ws.Range("S4").Value = 0.04361727704806999
$ws.Range("T4").Value = 0.04361727704807
$ws.Range("I5").Value = 0.852693657476708
$ws.Range("J5").Value = 0.852693657476708
$ws.Range("M5").Value = 1.363346333333333
$ws.Range("N5").Value = 4.090039
$ws.Range("O5").Value = 0.02430403345239443
$ws.Range("P5").Value = 0.02430403345239443
$ws.Range("Q5").Value = 0.2322092375323333
$ws.Range("R5").Value = 2.089883137791
$ws.Range("S5").Value = 0.02072389517595847
$ws.Range("T5").Value = 0.02072389517595847
$ws.Range("I6").Value = 0.852693657476708
$ws.Range("J6").Value = 0.852693657476708
$ws.Range("O6").Value = 0.679596855668023
$ws.Range("P6").Value = 0.679596855668023
$ws.Range("S6").Value = 0.5794879284692369
$ws.Range("T6").Value = 0.5794879284692369
$ws.Range("I7").Value = 0.852693657476708
$ws.Range("J7").Value = 0.852693657476708
$ws.Range("N7").Value = 49.82946200000001
$ws.Range("O7").Value = 0.2960991108795826
$ws.Range("P7").Value = 0.2960991108795826
$ws.Range("S7").Value = 0.2524818338315126
$ws.Range("T7").Value = 0.2524818338315126
